$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")

# --- Row 2 (hard coal): update years 2018-2050 (columns E:AK) ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 4
$ws.Range("G2").Value = 4
$ws.Range("H2").Formula = "=G2"
$ws.Range("I2:AK2").Formula = "=H2"

# Colour the updated cells (font colour -> theme accent, matches OOXML theme index 6)
$ws.Range("E2:AK2").Font.ThemeColor = 7

# --- Row 11 (petroleum): bump base priority from 1 to 2 ---
$ws.Range("B11").Value = 2

# --- Restore selections as left by the author ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("B22").Select()

$ws.Activate()
$ws.Range("E23").Select()
